$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data (TPM-updated) for rows 2-9, columns A-T.
# Column order:
# A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E Ligand-expressing cells, F Ligand detection rate, G Ligand average expression value,
# H Ligand total expression value, I Ligand derived specificity of average expression value,
# J Ligand derived specificity of total expression value, K Receptor-expressing cells,
# L Receptor detection rate, M Receptor average expression value, N Receptor total expression value,
# O Receptor derived specificity of average expression value, P Receptor derived specificity of total expression value,
# Q Edge average expression weight, R Edge total expression weight,
# S Edge average expression derived specificity, T Edge total expression derived specificity

$data = @(
    @("FAPs","a","Mc5r","ECs",2,0.6666666666666666,0.344207,1.032621,0.6985282229833164,0.6985282229833165,3,1,0.6692483333333333,2.007745,0.162783276239793,0.162783276239793,0.2303599610716666,2.073239649645,0.1137087126831849,0.113708712683185),
    @("FAPs","a","Mc5r","FAPs",2,0.6666666666666666,0.344207,1.032621,0.6985282229833164,0.6985282229833165,3,1,2.534005,7.602015,0.6163536244513372,0.6163536244513373,0.8722222590349999,7.850000331315,0.430540402017319,0.4305404020173191),
    @("FAPs","a","Mc5r","MuSCs",2,0.6666666666666666,0.344207,1.032621,0.6985282229833164,0.6985282229833165,3,1,0.885104,2.655312,0.2152864964419472,0.2152864964419472,0.304658992528,2.741930932752,0.1503836937918974,0.1503836937918975),
    @("FAPs","a","Mc5r","Resolving-Mac",2,0.6666666666666666,0.344207,1.032621,0.6985282229833164,0.6985282229833165,1,0.3333333333333333,0.022927,0.06878099999999999,0.005576602866922444,0.005576602866922445,0.007891633888999999,0.071024705001,0.003895414490915003,0.003895414490915004),
    @("MuSCs","a","Mc5r","ECs",1,0.3333333333333333,0.1485533333333333,0.44566,0.3014717770166836,0.3014717770166836,3,1,0.6692483333333333,2.007745,0.162783276239793,0.162783276239793,0.09941907074444445,0.8947716366999999,0.04907456355660809,0.04907456355660809),
    @("MuSCs","a","Mc5r","FAPs",1,0.3333333333333333,0.1485533333333333,0.44566,0.3014717770166836,0.3014717770166836,3,1,2.534005,7.602015,0.6163536244513372,0.6163536244513373,0.3764348894333334,3.3879140049,0.1858132224340183,0.1858132224340183),
    @("MuSCs","a","Mc5r","MuSCs",1,0.3333333333333333,0.1485533333333333,0.44566,0.3014717770166836,0.3014717770166836,3,1,0.885104,2.655312,0.2152864964419472,0.2152864964419472,0.1314851495466667,1.18336634592,0.06490280265004975,0.06490280265004976),
    @("MuSCs","a","Mc5r","Resolving-Mac",1,0.3333333333333333,0.1485533333333333,0.44566,0.3014717770166836,0.3014717770166836,1,0.3333333333333333,0.022927,0.06878099999999999,0.005576602866922444,0.005576602866922445,0.003405882273333334,0.03065294046,0.001681188376007442,0.001681188376007442)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $rowVals = $data[$i]
    for ($c = 1; $c -le $rowVals.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowVals[$c - 1]
    }
}
